$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.442.84'
$ws.Range("E2").Value = '  -2.54%  '

$ws.Range("D3").Value = '2.384.85'
$ws.Range("E3").Value = '  -4.72%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.46%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.571'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.34%  '

$ws.Range("D9").Value = '2.389.72'
$ws.Range("E9").Value = '  -4.39%  '

$ws.Range("E10").Value = '  -4.59%  '

$ws.Range("E11").Value = '  -0.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.92%  '

$ws.Range("E13").Value = '  -5.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.32%  '

$ws.Range("D15").Value = '2.816.91'
$ws.Range("E15").Value = '  -4.56%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000162'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.24%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '60.395.11'
$ws.Range("E17").Value = '  -2.38%  '

$ws.Range("D18").Value = '2.385.59'
$ws.Range("E18").Value = '  -4.63%  '

$ws.Range("E19").Value = '  -5.85%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.52%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '311.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.28%  '

$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.16%  '

$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("E24").Value = '  +0.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.50'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.74%  '

$ws.Range("D27").Value = '2.509.59'
$ws.Range("E27").Value = '  -4.33%  '

$ws.Range("D28").Value = '0.0₃0900'
$ws.Range("E28").Value = '  -11.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.03%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.52%  '

$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.40'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '498.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.23%  '

$ws.Range("E33").Value = '  -4.82%  '

$ws.Range("E34").Value = '  -4.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.25%  '

$ws.Range("E36").Value = '  +0.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.83%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.57'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.370'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '136.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.49%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.48%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.25%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.18%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '139.08'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.49'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.87%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.94'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0510'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.574'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0917'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.86%  '
